$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Bengali character table (columns A, D, G, J, M) with the new
# shifted set of letters, and fill in the newly added row values for
# columns M/N on rows 6-8, plus three new letters in column A rows 6-8.

$ws.Range("A1").Value = "প"
$ws.Range("D1").Value = "ক"
$ws.Range("G1").Value = "ঝ"
$ws.Range("J1").Value = "থ"
$ws.Range("M1").Value = "উ"

$ws.Range("A2").Value = "ফ"
$ws.Range("D2").Value = "খ"
$ws.Range("G2").Value = "ঞ"
$ws.Range("J2").Value = "দ"
$ws.Range("M2").Value = "ঊ"

$ws.Range("A3").Value = "ব"
$ws.Range("D3").Value = "গ"
$ws.Range("G3").Value = "ট"
$ws.Range("J3").Value = "ধ"
$ws.Range("M3").Value = "ঋ"

$ws.Range("A4").Value = "ভ"
$ws.Range("D4").Value = "ঘ"
$ws.Range("G4").Value = "ঠ"
$ws.Range("J4").Value = "ন"
$ws.Range("M4").Value = "এ"

$ws.Range("A5").Value = "ম"
$ws.Range("D5").Value = "ঙ"
$ws.Range("G5").Value = "ড"
$ws.Range("J5").Value = "অ"
$ws.Range("M5").Value = "ঐ"

$ws.Range("A6").Value = "য"
$ws.Range("D6").Value = "চ"
$ws.Range("G6").Value = "ঢ"
$ws.Range("J6").Value = "আ"
$ws.Range("M6").Value = "ও"
$ws.Range("N6").Value = 37

$ws.Range("A7").Value = "র"
$ws.Range("D7").Value = "ছ"
$ws.Range("G7").Value = "ণ"
$ws.Range("J7").Value = "ই"
$ws.Range("M7").Value = "ঔ"
$ws.Range("N7").Value = 38

$ws.Range("A8").Value = "ল"
$ws.Range("D8").Value = "জ"
$ws.Range("G8").Value = "ত"
$ws.Range("J8").Value = "ঈ"
$ws.Range("M8").Value = "ঌ"
$ws.Range("N8").Value = 39

# Update the view: scroll so row 2 is the top visible row, and move the
# active selection to N8.
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("N8").Select()
